# process updated for USA process
# Adds new "mail notification" config rows to the Constants sheet, wires up a
# mailto hyperlink on the MailTo value, and leaves the UI focused on the
# Constants sheet/cell B8 (matching the saved view state of the source file).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# --- new label / value rows -------------------------------------------------
# (values are written in this specific order so the shared-strings table ends
# up in the same sequence as the authored workbook)

$ws.Range("A6").Value = "DataNoBelongToCurrentWeekMailSubject"
$ws.Range("B6").Value = "Notification- Data doesn't belong to this week"

$ws.Range("A7").Value = "DataNoBelongToCurrentWeekMailBody"

$ws.Range("A9").Value = "MailTo"
$ws.Range("B9").Value = "manindersinghbisht77@gmail.com"

$ws.Range("A10").Value = "MailCC"

$ws.Range("A12").Value = "InputFileExceptionSubject"
$ws.Range("B12").Value = "Error- Input file Notification"

$ws.Range("A13").Value = "InputFileExceptionbody"

$ws.Range("B15").Value = "Notification- Not All file Got downloaded"
$ws.Range("A15").Value = "downloadFileExceptionSubject"

$ws.Range("A16").Value = "downloadFileExceptionbody"

$bodyNoBelong = @"
Dear Team<br/>
Please find the below data which does not belong to current week.<br/>
[Nothisweekdatatable]<br>
Also let us know if anything is required<br/>
Thank you,<br/>
**********************This is system generated E-Mail, please do not respond on this************
"@
$ws.Range("B7").Value = $bodyNoBelong

$bodyInputException = @"
Dear Team<br/>
There is a error in Cinepolish USA run, below is the detail of error<br/>
[error]<br/>
Also let us know if anything is required<br/>
Thank you,<br/>
**********************This is system generated E-Mail, please do not respond on this************
"@
$ws.Range("B13").Value = $bodyInputException

$bodyDownloadException = @"
Dear Team<br/>
There is a error in Cinepos USA process run, below is the detail of error<br/>
Not all required files got downloaded from FTP<br/>
Also let us know if anything is required<br/>
Thank you,<br/>
**********************This is system generated E-Mail, please do not respond on this************
"@
$ws.Range("B16").Value = $bodyDownloadException

# --- formatting --------------------------------------------------------------
# B7/B13/B16 hold the long multi-line mail bodies -> wrap text, but keep the
# original row height instead of letting it auto-fit.
$ws.Range("B7").WrapText = $true
$ws.Range("B13").WrapText = $true
$ws.Range("B16").WrapText = $true
$ws.Rows.Item(7).RowHeight = 14.25
$ws.Rows.Item(13).RowHeight = 14.25
$ws.Rows.Item(16).RowHeight = 14.25

# B8 stays empty but keeps the underlined "placeholder" look it already had.
$ws.Range("A8").Font.Underline = $false
$ws.Range("B8").Font.Underline = $true

# MailTo value becomes a live mailto: link.
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:manindersinghbisht77@gmail.com") | Out-Null

# --- view state ---------------------------------------------------------------
# Settings keeps its original selection (B2) but the saved view had scrolled
# down; Constants becomes the active/visible tab with B8 selected.
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Activate()
$wsSettings.Range("B2").Select() | Out-Null

$ws.Activate()
$ws.Range("B8").Select() | Out-Null
